$p = $ppt.ActivePresentation

# --- Slide 1: remove the two small flower/leaf decoration groups (602, 605) ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).Delete()
$s1.Shapes.Item(2).Delete()

# Add three new rectangles ("Правоъгълник 1/2/3") at the end of the shape tree
$rect1 = $s1.Shapes.AddShape(1, 142.0897, 119.0769, 59.2616, 44.8615)
$rect1.Name = "Правоъгълник 1"
$rect1.Fill.ForeColor.RGB = 14479093
$rect1.Line.Visible = $false
$rect1.TextFrame.VerticalAnchor = 3
$rect1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$rect2 = $s1.Shapes.AddShape(1, 575.5665, 185.2999, 59.2616, 44.8615)
$rect2.Name = "Правоъгълник 2"
$rect2.Fill.ForeColor.RGB = 14479093
$rect2.Line.Visible = $false
$rect2.TextFrame.VerticalAnchor = 3
$rect2.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$rect3 = $s1.Shapes.AddShape(1, 575.5665, 212.4923, 59.2616, 44.8615)
$rect3.Name = "Правоъгълник 3"
$rect3.Fill.ForeColor.RGB = 14479093
$rect3.Line.Visible = $false
$rect3.TextFrame.VerticalAnchor = 3
$rect3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Slide 4: remove the three small flower/leaf decoration groups (646, 649, 652) ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(16).Delete()
$s4.Shapes.Item(15).Delete()
$s4.Shapes.Item(14).Delete()

# --- Slide 5: remove the trailing decoration group (id 50, "Google Shape;646;p29") ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item($s5.Shapes.Count).Delete()

# --- Slide 6: remove the trailing decoration group (id 22, "Google Shape;774;p34") ---
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item($s6.Shapes.Count).Delete()
